$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Trf", "Tfrc", "ECs", 2, 0.6666666666666666, 14.52107966666667, 43.563239, 0.2237994865827619, 0.2237994865827619, 3, 1, 27.646565, 82.939695, 0.4490272892980421, 0.4490272892980421, 401.4579728746783, 3613.121755872105, 0.1004920768065511, 0.1004920768065511),
    @("ECs", "Trf", "Tfrc", "FAPs", 2, 0.6666666666666666, 14.52107966666667, 43.563239, 0.2237994865827619, 0.2237994865827619, 3, 1, 3.745931333333334, 11.237794, 0.06084030303595646, 0.06084030303595646, 54.39496731719623, 489.5547058547661, 0.0136160285829867, 0.0136160285829867),
    @("ECs", "Trf", "Tfrc", "M2", 2, 0.6666666666666666, 14.52107966666667, 43.563239, 0.2237994865827619, 0.2237994865827619, 3, 1, 15.12837333333333, 45.38512, 0.2457105419554094, 0.2457105419554094, 219.6803144004089, 1977.12282960368, 0.0549898931375928, 0.0549898931375928),
    @("ECs", "Trf", "Tfrc", "sCs", 2, 0.6666666666666666, 14.52107966666667, 43.563239, 0.2237994865827619, 0.2237994865827619, 3, 1, 15.04902966666666, 45.14708899999999, 0.2444218657105919, 0.2444218657105919, 218.5281586956968, 1966.753428261271, 0.05470148805563125, 0.05470148805563125),
    @("FAPs", "Trf", "Tfrc", "ECs", 3, 1, 9.234235, 27.702705, 0.1423184156704627, 0.1423184156704627, 3, 1, 27.646565, 82.939695, 0.4490272892980421, 0.4490272892980421, 255.294878152775, 2297.653903374975, 0.06390485240569985, 0.06390485240569985),
    @("FAPs", "Trf", "Tfrc", "FAPs", 3, 1, 9.234235, 27.702705, 0.1423184156704627, 0.1423184156704627, 3, 1, 3.745931333333334, 11.237794, 0.06084030303595646, 0.06084030303595646, 34.59081022586334, 311.31729203277, 0.008658695536988163, 0.008658695536988163),
    @("FAPs", "Trf", "Tfrc", "M2", 3, 1, 9.234235, 27.702705, 0.1423184156704627, 0.1423184156704627, 3, 1, 15.12837333333333, 45.38512, 0.2457105419554094, 0.2457105419554094, 139.6989545277333, 1257.2905907496, 0.03496913504462461, 0.03496913504462461),
    @("FAPs", "Trf", "Tfrc", "sCs", 3, 1, 9.234235, 27.702705, 0.1423184156704627, 0.1423184156704627, 3, 1, 15.04902966666666, 45.14708899999999, 0.2444218657105919, 0.2444218657105919, 138.9662764639716, 1250.696488175745, 0.03478573268315003, 0.03478573268315003),
    @("M2", "Trf", "Tfrc", "ECs", 3, 1, 40.17359466666667, 120.520784, 0.6191571196474152, 0.6191571196474152, 3, 1, 27.646565, 82.939695, 0.4490272892980421, 0.4490272892980421, 1110.661896235653, 9995.95706612088, 0.2780184430848624, 0.2780184430848624),
    @("M2", "Trf", "Tfrc", "FAPs", 3, 1, 40.17359466666667, 120.520784, 0.6191571196474152, 0.6191571196474152, 3, 1, 3.745931333333334, 11.237794, 0.06084030303595646, 0.06084030303595646, 150.4875270344995, 1354.387743310496, 0.03766970678621869, 0.03766970678621869),
    @("M2", "Trf", "Tfrc", "M2", 3, 1, 40.17359466666667, 120.520784, 0.6191571196474152, 0.6191571196474152, 3, 1, 15.12837333333333, 45.38512, 0.2457105419554094, 0.2457105419554094, 607.7611382593423, 5469.85024433408, 0.1521334314241166, 0.1521334314241166),
    @("M2", "Trf", "Tfrc", "sCs", 3, 1, 40.17359466666667, 120.520784, 0.6191571196474152, 0.6191571196474152, 3, 1, 15.04902966666666, 45.14708899999999, 0.2444218657105919, 0.2444218657105919, 604.5736179553084, 5441.162561597775, 0.1513355383522174, 0.1513355383522174),
    @("sCs", "Trf", "Tfrc", "ECs", 3, 1, 0.9554203333333334, 2.866261, 0.01472497809936019, 0.01472497809936019, 3, 1, 27.646565, 82.939695, 0.4490272892980421, 0.4490272892980421, 26.41409034782167, 237.726813130395, 0.006611917000928742, 0.006611917000928742),
    @("sCs", "Trf", "Tfrc", "FAPs", 3, 1, 0.9554203333333334, 2.866261, 0.01472497809936019, 0.01472497809936019, 3, 1, 3.745931333333334, 11.237794, 0.06084030303595646, 0.06084030303595646, 3.578938963137112, 32.210450668234, 0.0008958721297628961, 0.0008958721297628961),
    @("sCs", "Trf", "Tfrc", "M2", 3, 1, 0.9554203333333334, 2.866261, 0.01472497809936019, 0.01472497809936019, 3, 1, 15.12837333333333, 45.38512, 0.2457105419554094, 0.2457105419554094, 14.45395549292445, 130.08559943632, 0.003618082349075327, 0.003618082349075326),
    @("sCs", "Trf", "Tfrc", "sCs", 3, 1, 0.9554203333333334, 2.866261, 0.01472497809936019, 0.01472497809936019, 3, 1, 15.04902966666666, 45.14708899999999, 0.2444218657105919, 0.2444218657105919, 14.37814894046989, 129.403340464229, 0.003599106619593223, 0.003599106619593223)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}

Write-Host "Applied $($data.Count) rows"